$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.442.02'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '3.504.29'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '585.22'
$ws.Range("E5").Value = '  -2.35%  '
$ws.Range("D6").Value = '175.21'
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.498.87'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.595'
$ws.Range("E9").Value = '  -3.14%  '
$ws.Range("E10").Value = '  -3.66%  '
$ws.Range("D11").Value = '6.87'
$ws.Range("E11").Value = '  -2.05%  '
$ws.Range("D12").Value = '0.421'
$ws.Range("E12").Value = '  -3.68%  '
$ws.Range("D13").Value = '4.110.78'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").Value = '30.42'
$ws.Range("E14").Value = '  -6.19%  '
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("D16").Value = '66.431.81'
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("D17").Value = '0.0000173'
$ws.Range("E17").Value = '  -3.06%  '
$ws.Range("D18").Value = '3.500.49'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '6.01'
$ws.Range("E19").Value = '  -5.08%  '
$ws.Range("D20").Value = '13.88'
$ws.Range("E20").Value = '  -3.18%  '
$ws.Range("D21").Value = '379.59'
$ws.Range("E21").Value = '  -3.21%  '
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").Value = '0.548'
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = '5.76'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").Value = '72.08'
$ws.Range("E26").Value = '  -1.51%  '
$ws.Range("D27").Value = '0.0000121'
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("D28").Value = '9.87'
$ws.Range("E28").Value = '  -5.03%  '
$ws.Range("D29").Value = '0.173'
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '24.45'
$ws.Range("E31").Value = '  +3.60%  '
$ws.Range("D32").Value = '5.92'
$ws.Range("E32").Value = '  -3.96%  '
$ws.Range("E33").Value = '  -3.18%  '
$ws.Range("D34").Value = '1.33'
$ws.Range("E34").Value = '  -6.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '7.24'
$ws.Range("E36").Value = '  -2.83%  '
$ws.Range("D37").Value = '1.57'
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("D38").Value = '29.63'
$ws.Range("E38").Value = '  +11.97%  '
$ws.Range("D39").Value = '160.09'
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("D41").Value = '1.79'
$ws.Range("E41").Value = '  -5.74%  '
$ws.Range("E42").Value = '  -2.60%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.55'
$ws.Range("E43").Value = '  -10.85%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '6.43'
$ws.Range("E44").Value = '  -6.18%  '
$ws.Range("D45").Value = '2.694.67'
$ws.Range("E45").Value = '  -5.37%  '
$ws.Range("D46").Value = '0.0697'
$ws.Range("E46").Value = '  -4.02%  '
$ws.Range("D47").Value = '40.69'
$ws.Range("E47").Value = '  -2.33%  '
$ws.Range("D48").Value = '24.73'
$ws.Range("E48").Value = '  -7.90%  '
$ws.Range("D49").Value = '0.0292'
$ws.Range("E49").Value = '  -3.14%  '
$ws.Range("D50").Value = '317.43'
$ws.Range("E50").Value = '  -5.60%  '
$ws.Range("E51").Value = '  -5.31%  '
